$d = $word.ActiveDocument

# Locate the paragraph that contains the sentence we need to split into
# several runs (Find on the exact original run text so we don't
# accidentally match anything else in the document).
$r = $d.Content
$found = $r.Find.Execute(" The only thing I got from this package was models and particle effects.")
if (-not $found) {
    throw "Could not find target sentence"
}

# Grab the whole containing paragraph (the Find range is narrowed to the
# matched text, so Paragraphs(1) gives us the single paragraph it lives in).
$para = $r.Paragraphs(1)
$pr = $para.Range

$ooxml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w14:paraId="72287480" w14:textId="002289F0" w:rsidR="00B05F4F" w:rsidRPr="00B05F4F" w:rsidRDefault="00B05F4F" w:rsidP="00B05F4F"><w:r><w:t>This asset package has many resources including scripts  that govern player movement, gun firing ,visual effects, lighting and much more. However I only used the models for the handgun and the rifle and the particle system for the muzzle flash. All the code that allow the assets to function in my prototype was written by me, animations created by me and lighting for the muzzle flash created by me.</w:t></w:r><w:r w:rsidR="0064382A"><w:t xml:space="preserve"> The only thing I got from this package was models</w:t></w:r><w:r><w:t>, the sound effects</w:t></w:r><w:r><w:t xml:space="preserve"> and</w:t></w:r><w:r><w:t xml:space="preserve"> some</w:t></w:r><w:r><w:t xml:space="preserve"> particle effects.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$pr.InsertXML($ooxml)
